$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "# of Games" counts in column B (rows 2-31) ---
$updates = @{
    3  = 4
    4  = 4
    5  = 4
    6  = 4
    7  = 4
    9  = 7
    10 = 7
    11 = 5
    12 = 2
    13 = 3
    14 = 5
    17 = 2
    18 = 4
    19 = 3
    20 = 5
    21 = 6
    22 = 2
    23 = 3
    24 = 4
    27 = 3
    28 = 1
    29 = 1
    31 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# --- Row 33 total stays a formula; recalculated automatically ---
$ws.Range("B33").Formula = "=SUM(B2:B31)/2"

# --- New summary rows ---
$ws.Range("A34").Value = "Average per Team"
$ws.Range("B34").Formula = "=AVERAGE(B2:B31)"
$ws.Range("B34").NumberFormat = "0.0"

$ws.Range("A35").Value = "Max Team"
$ws.Range("B35").Formula = "=MAX(B2:B31)"

$ws.Range("A36").Value = "Min Team"
$ws.Range("B36").Formula = "=MIN(B2:B31)"

# --- Column A width adjustment (now fits the longer "Average per Team" label) ---
$ws.Columns.Item(1).ColumnWidth = 15.75

# --- Update the active selection to B15 ---
$ws.Range("B15").Select()
